# Update the cryptos list (Price column D and Volume(1h) column E)
# with refreshed values, matching a new GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @(
    @{ Cell = 'D2'; Value = '''20.561.39'; ForceText = $true },
    @{ Cell = 'E2'; Value = '  +1.79%  '; ForceText = $false },
    @{ Cell = 'D3'; Value = '''1.473.64'; ForceText = $true },
    @{ Cell = 'E3'; Value = '  +2.69%  '; ForceText = $false },
    @{ Cell = 'E4'; Value = '  -0.17%  '; ForceText = $false },
    @{ Cell = 'D5'; Value = '''0.9576'; ForceText = $true },
    @{ Cell = 'E5'; Value = '  +4.88%  '; ForceText = $false },
    @{ Cell = 'D6'; Value = '''277.65'; ForceText = $true },
    @{ Cell = 'E6'; Value = '  +0.37%  '; ForceText = $false },
    @{ Cell = 'D7'; Value = '''0.3627'; ForceText = $true },
    @{ Cell = 'E7'; Value = '  -0.28%  '; ForceText = $false },
    @{ Cell = 'D8'; Value = '''0.3083'; ForceText = $true },
    @{ Cell = 'E8'; Value = '  -0.34%  '; ForceText = $false },
    @{ Cell = 'D9'; Value = '''39.65'; ForceText = $true },
    @{ Cell = 'E9'; Value = '  +1.66%  '; ForceText = $false },
    @{ Cell = 'D10'; Value = '''1.074'; ForceText = $true },
    @{ Cell = 'E10'; Value = '  +5.70%  '; ForceText = $false },
    @{ Cell = 'E11'; Value = '  +2.16%  '; ForceText = $false },
    @{ Cell = 'E12'; Value = '  -0.14%  '; ForceText = $false },
    @{ Cell = 'D13'; Value = '''5.505'; ForceText = $true },
    @{ Cell = 'E13'; Value = '  +3.11%  '; ForceText = $false },
    @{ Cell = 'E14'; Value = '  +4.37%  '; ForceText = $false },
    @{ Cell = 'D15'; Value = '''0.9589'; ForceText = $true },
    @{ Cell = 'E15'; Value = '  +1.90%  '; ForceText = $false },
    @{ Cell = 'D16'; Value = '''6.163'; ForceText = $true },
    @{ Cell = 'E16'; Value = '  +2.02%  '; ForceText = $false },
    @{ Cell = 'E17'; Value = '  +1.44%  '; ForceText = $false },
    @{ Cell = 'D18'; Value = '''1.472.19'; ForceText = $true },
    @{ Cell = 'E18'; Value = '  +2.37%  '; ForceText = $false },
    @{ Cell = 'D19'; Value = '''0.05941'; ForceText = $true },
    @{ Cell = 'E19'; Value = '  +5.37%  '; ForceText = $false },
    @{ Cell = 'D20'; Value = '''68.95'; ForceText = $true },
    @{ Cell = 'E20'; Value = '  +1.83%  '; ForceText = $false },
    @{ Cell = 'D21'; Value = '''5.503'; ForceText = $true },
    @{ Cell = 'E21'; Value = '  +2.91%  '; ForceText = $false },
    @{ Cell = 'E22'; Value = '  +1.97%  '; ForceText = $false },
    @{ Cell = 'D23'; Value = '''11.19'; ForceText = $true },
    @{ Cell = 'E23'; Value = '  +3.93%  '; ForceText = $false },
    @{ Cell = 'D24'; Value = '''2.265'; ForceText = $true },
    @{ Cell = 'E24'; Value = '  +0.74%  '; ForceText = $false },
    @{ Cell = 'D25'; Value = '''20.576.83'; ForceText = $true },
    @{ Cell = 'E25'; Value = '  +1.51%  '; ForceText = $false },
    @{ Cell = 'D26'; Value = '''142.31'; ForceText = $true },
    @{ Cell = 'E26'; Value = '  +4.23%  '; ForceText = $false },
    @{ Cell = 'D27'; Value = '''2.130'; ForceText = $true },
    @{ Cell = 'E27'; Value = '  -0.19%  '; ForceText = $false },
    @{ Cell = 'D28'; Value = '''17.17'; ForceText = $true },
    @{ Cell = 'E28'; Value = '  +1.80%  '; ForceText = $false },
    @{ Cell = 'D29'; Value = '''1.634.13'; ForceText = $true },
    @{ Cell = 'E29'; Value = '  +2.81%  '; ForceText = $false },
    @{ Cell = 'D30'; Value = '''113.84'; ForceText = $true },
    @{ Cell = 'E30'; Value = '  +3.97%  '; ForceText = $false },
    @{ Cell = 'D31'; Value = '''3.904'; ForceText = $true },
    @{ Cell = 'E31'; Value = '  +0.09%  '; ForceText = $false },
    @{ Cell = 'D32'; Value = '''0.08023'; ForceText = $true },
    @{ Cell = 'E32'; Value = '  +4.76%  '; ForceText = $false },
    @{ Cell = 'D33'; Value = '''4.941'; ForceText = $true },
    @{ Cell = 'E33'; Value = '  +2.82%  '; ForceText = $false },
    @{ Cell = 'D34'; Value = '''0.8021'; ForceText = $true },
    @{ Cell = 'E34'; Value = '  +1.30%  '; ForceText = $false },
    @{ Cell = 'D35'; Value = '''1.511'; ForceText = $true },
    @{ Cell = 'E35'; Value = '  +4.54%  '; ForceText = $false },
    @{ Cell = 'D36'; Value = '''1.213'; ForceText = $true },
    @{ Cell = 'E36'; Value = '  +6.73%  '; ForceText = $false },
    @{ Cell = 'D37'; Value = '''0.05762'; ForceText = $true },
    @{ Cell = 'E37'; Value = '  -2.11%  '; ForceText = $false },
    @{ Cell = 'D38'; Value = '''4.722'; ForceText = $true },
    @{ Cell = 'E38'; Value = '  +2.64%  '; ForceText = $false },
    @{ Cell = 'E39'; Value = '  +3.72%  '; ForceText = $false },
    @{ Cell = 'D40'; Value = '''0.9590'; ForceText = $true },
    @{ Cell = 'E40'; Value = '  +4.05%  '; ForceText = $false },
    @{ Cell = 'D41'; Value = '''10.40'; ForceText = $true },
    @{ Cell = 'E41'; Value = '  +2.72%  '; ForceText = $false },
    @{ Cell = 'D42'; Value = '''0.1883'; ForceText = $true },
    @{ Cell = 'E42'; Value = '  +2.72%  '; ForceText = $false },
    @{ Cell = 'D43'; Value = '''7.425'; ForceText = $true },
    @{ Cell = 'E43'; Value = '  +6.10%  '; ForceText = $false },
    @{ Cell = 'D44'; Value = '''0.5290'; ForceText = $true },
    @{ Cell = 'E44'; Value = '  +1.77%  '; ForceText = $false },
    @{ Cell = 'D45'; Value = '''3.526'; ForceText = $true },
    @{ Cell = 'E45'; Value = '  +0.73%  '; ForceText = $false },
    @{ Cell = 'D46'; Value = '''12.23'; ForceText = $true },
    @{ Cell = 'E46'; Value = '  +2.64%  '; ForceText = $false },
    @{ Cell = 'D47'; Value = '''118.71'; ForceText = $true },
    @{ Cell = 'E47'; Value = '  +0.99%  '; ForceText = $false },
    @{ Cell = 'D48'; Value = '''0.5211'; ForceText = $true },
    @{ Cell = 'E48'; Value = '  +2.40%  '; ForceText = $false },
    @{ Cell = 'E49'; Value = '  +3.99%  '; ForceText = $false },
    @{ Cell = 'D50'; Value = '''0.06462'; ForceText = $true },
    @{ Cell = 'E50'; Value = '  +2.35%  '; ForceText = $false },
    @{ Cell = 'D51'; Value = '''0.9879'; ForceText = $true },
    @{ Cell = 'E51'; Value = '  +0.11%  '; ForceText = $false }
)

foreach ($u in $updates) {
    $ws.Range($u.Cell).Value = $u.Value
    if ($u.ForceText) {
        # The leading apostrophe above forces Excel to keep the numeric-looking
        # text (e.g. "1.074") as a string instead of converting it to a number.
        # Resetting the style back to Normal drops the transient "quote prefix"
        # formatting flag that the apostrophe entry leaves behind, so the cell
        # ends up with the same (default) styling as before the edit.
        $ws.Range($u.Cell).Style = "Normal"
    }
}
